# "Generate Report for Handoff"
#
# The localization-status report is refreshed: the cb4dc240-....md file's
# handoff/status information moves from "Handed back: in sync with en-US"
# to "Ready for handoff" (a new handoff round was generated), the relevant
# "Latest Handoff Datetime" timestamps are updated, and the zh-cn / de-de
# detail sheets record an Error Detail note about the handback file being
# stale. The Error Detail column is also widened so the new long message is
# readable.

$wb = $excel.ActiveWorkbook

$longError = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/621ae638fdf405cea936e34bc356a9f494abfea4/e2e/cb4dc240-bbcf-4a5b-8475-8f268fea9a70.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20de6d281dc5a563e47e093993cd98c9e00ecd22/e2e/cb4dc240-bbcf-4a5b-8475-8f268fea9a70.md."

# --- Overview sheet: cb4dc240 row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-30 06:51:16"

# --- zh-cn sheet: cb4dc240 row (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-30 06:51:11"
$wsZhCn.Range("P3").Value = $longError
# ColumnWidth is entered in "characters"; Excel then rounds it to a whole
# number of pixels before persisting it to xlsx, so 39.17 characters is what
# lands on a saved width of exactly 40 (matching the target column width).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: cb4dc240 row (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-30 06:51:16"
$wsDeDe.Range("P3").Value = $longError
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
